# Update crypto price (D) and volume-change (E) columns to the latest scrape.
# Values that look like plain decimals (e.g. "226.00") need a leading quote
# so Excel stores them as text (matching the original inline-string cells)
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.524.85"
$ws.Range("E2").Value = "  +5.24%  "

$ws.Range("D3").Value = "1.724.66"
$ws.Range("E3").Value = "  +4.10%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'226.00"
$ws.Range("E5").Value = "  +3.32%  "

$ws.Range("D6").Value = "'0.5374"
$ws.Range("E6").Value = "  +2.59%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.2675"

$ws.Range("D9").Value = "'0.06602"
$ws.Range("E9").Value = "  +3.83%  "

$ws.Range("D10").Value = "'21.75"
$ws.Range("E10").Value = "  +5.60%  "

$ws.Range("D11").Value = "'0.07743"

$ws.Range("D12").Value = "'4.620"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").Value = "1.742.90"
$ws.Range("E13").Value = "  +4.69%  "

$ws.Range("E14").Value = "  +4.13%  "

$ws.Range("D15").Value = "'0.5869"
$ws.Range("E15").Value = "  +4.24%  "

$ws.Range("D16").Value = "0.0₅8314"
$ws.Range("E16").Value = "  +1.28%  "

$ws.Range("D17").Value = "'68.07"
$ws.Range("E17").Value = "  +3.93%  "

$ws.Range("D18").Value = "27.537.95"
$ws.Range("E18").Value = "  +5.31%  "

$ws.Range("D19").Value = "'222.64"
$ws.Range("E19").Value = "  +15.55%  "

$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "'4.746"
$ws.Range("E21").Value = "  +1.75%  "

$ws.Range("D22").Value = "'10.69"
$ws.Range("E22").Value = "  +1.41%  "

$ws.Range("D23").Value = "'6.097"
$ws.Range("E23").Value = "  +2.27%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "'148.10"
$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("E26").Value = "  +12.32%  "

$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("D28").Value = "'7.399"
$ws.Range("E28").Value = "  +1.78%  "

$ws.Range("D29").Value = "'16.70"
$ws.Range("E29").Value = "  +4.50%  "

$ws.Range("D30").Value = "'0.05535"
$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("E31").Value = "  +2.42%  "

$ws.Range("E32").Value = "  +2.22%  "

$ws.Range("D33").Value = "'3.464"
$ws.Range("E33").Value = "  +2.42%  "

$ws.Range("D34").Value = "'1.661"
$ws.Range("E34").Value = "  +6.01%  "

$ws.Range("D35").Value = "'0.9598"
$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("E36").Value = "  +1.53%  "

$ws.Range("E37").Value = "  +1.83%  "

$ws.Range("D38").Value = "'0.5943"
$ws.Range("E38").Value = "  +4.38%  "

$ws.Range("E39").Value = "  +3.75%  "

$ws.Range("D40").Value = "'5.918"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").Value = "1.060.44"
$ws.Range("E41").Value = "  +3.17%  "

$ws.Range("D42").Value = "'0.8559"
$ws.Range("E42").Value = "  +2.64%  "

$ws.Range("D43").Value = "'1.004"
$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").Value = "'101.47"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").Value = "1.868.13"
$ws.Range("E45").Value = "  +4.04%  "

$ws.Range("E46").Value = "  +10.38%  "

$ws.Range("D47").Value = "'59.08"
$ws.Range("E47").Value = "  +2.16%  "

$ws.Range("D48").Value = "'8.218"
$ws.Range("E48").Value = "  +2.11%  "

$ws.Range("D49").Value = "'0.4441"
$ws.Range("E49").Value = "  +2.22%  "

$ws.Range("D50").Value = "'1.006"
$ws.Range("E50").Value = "  +0.32%  "

$ws.Range("D51").Value = "'0.05275"
$ws.Range("E51").Value = "  +1.56%  "
